$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-06-26 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-27 Thursday", 2) | Out-Null
$d.Content.Find.Execute("174×9=1566", $true, $false, $false, $false, $false, $true, 1, $false, "254×8=2032", 2) | Out-Null
$d.Content.Find.Execute("615×8=4920", $true, $false, $false, $false, $false, $true, 1, $false, "407×5=2035", 2) | Out-Null
$d.Content.Find.Execute("942×6=5652", $true, $false, $false, $false, $false, $true, 1, $false, "902×8=7216", 2) | Out-Null
$d.Content.Find.Execute("270×2=540", $true, $false, $false, $false, $false, $true, 1, $false, "164×2=328", 2) | Out-Null
$d.Content.Find.Execute("967×2=1934", $true, $false, $false, $false, $false, $true, 1, $false, "948×2=1896", 2) | Out-Null
$d.Content.Find.Execute("844×8=6752", $true, $false, $false, $false, $false, $true, 1, $false, "364×8=2912", 2) | Out-Null
$d.Content.Find.Execute("193×2=386", $true, $false, $false, $false, $false, $true, 1, $false, "899×3=2697", 2) | Out-Null
$d.Content.Find.Execute("882×3=2646", $true, $false, $false, $false, $false, $true, 1, $false, "337×2=674", 2) | Out-Null
$d.Content.Find.Execute("447×3=1341", $true, $false, $false, $false, $false, $true, 1, $false, "622×6=3732", 2) | Out-Null
$d.Content.Find.Execute("631×3=1893", $true, $false, $false, $false, $false, $true, 1, $false, "950×2=1900", 2) | Out-Null
$d.Content.Find.Execute("470×5=2350", $true, $false, $false, $false, $false, $true, 1, $false, "103×4=412", 2) | Out-Null
$d.Content.Find.Execute("493×2=986", $true, $false, $false, $false, $false, $true, 1, $false, "318×5=1590", 2) | Out-Null
$d.Content.Find.Execute("873×5=4365", $true, $false, $false, $false, $false, $true, 1, $false, "996×6=5976", 2) | Out-Null
$d.Content.Find.Execute("431×3=1293", $true, $false, $false, $false, $false, $true, 1, $false, "889×9=8001", 2) | Out-Null
$d.Content.Find.Execute("298×8=2384", $true, $false, $false, $false, $false, $true, 1, $false, "705×7=4935", 2) | Out-Null
$d.Content.Find.Execute("663×6=3978", $true, $false, $false, $false, $false, $true, 1, $false, "549×6=3294", 2) | Out-Null
$d.Content.Find.Execute("385×7=2695", $true, $false, $false, $false, $false, $true, 1, $false, "541×3=1623", 2) | Out-Null
$d.Content.Find.Execute("217×6=1302", $true, $false, $false, $false, $false, $true, 1, $false, "139×2=278", 2) | Out-Null
$d.Content.Find.Execute("174×7=1218", $true, $false, $false, $false, $false, $true, 1, $false, "143×9=1287", 2) | Out-Null
$d.Content.Find.Execute("793×4=3172", $true, $false, $false, $false, $false, $true, 1, $false, "435×4=1740", 2) | Out-Null
$d.Content.Find.Execute("372×2=744", $true, $false, $false, $false, $false, $true, 1, $false, "424×8=3392", 2) | Out-Null
$d.Content.Find.Execute("339×8=2712", $true, $false, $false, $false, $false, $true, 1, $false, "181×9=1629", 2) | Out-Null
$d.Content.Find.Execute("985×2=1970", $true, $false, $false, $false, $false, $true, 1, $false, "842×8=6736", 2) | Out-Null
$d.Content.Find.Execute("637×6=3822", $true, $false, $false, $false, $false, $true, 1, $false, "452×7=3164", 2) | Out-Null
$d.Content.Find.Execute("883×3=2649", $true, $false, $false, $false, $false, $true, 1, $false, "445×8=3560", 2) | Out-Null
